$wb = $excel.ActiveWorkbook

# --- Add the new "EXT - updated results" sheet as the very first sheet ---
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "EXT - updated results"

# Header row
$newSheet.Range("A1").Value = "Ontology"
$newSheet.Range("B1").Value = "CRF"
$newSheet.Range("C1").Value = "BIOBERT"
$newSheet.Range("D1").Value = "OGER"
$newSheet.Range("A1:D1").Font.Bold = $true

# Data rows: Ontology, CRF, BIOBERT, OGER
$data = @(
    @("CHEBI_EXT",      0.78910000000000002, 0.80389999999999995, 0.82089999999999996),
    @("CL_EXT",         0.73809999999999998, 0.74909999999999999, 0.74839999999999995),
    @("GO_BP_EXT",      0.72789999999999999, 0.73529999999999995, 0.81379999999999997),
    @("GO_CC_EXT",      0.87380000000000002, 0.89829999999999999, 0.89359999999999995),
    @("GO_MF_EXT",      0.64129999999999998, 0.62549999999999994, 0.74380000000000002),
    @("MOP_EXT",        0.8,                 0.86509999999999998, 0.84370000000000001),
    @("NCBITaxon_EXT",  0.871,               0.86240000000000006, 0.97219999999999995),
    @("PR_EXT",         0.43969999999999998, 0.51880000000000004, 0.80110000000000003),
    @("SO_EXT",         0.76819999999999999, 0.78290000000000004, 0.91869999999999996),
    @("UBERON_EXT",     0.75580000000000003, 0.77110000000000001, 0.77139999999999997)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("A$r").Font.Bold = $true
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("D$r").Font.Color = 0
    $r = $r + 1
}

$newSheet.Range("D11").Select()
